$d = $word.ActiveDocument

function Get-ParaByText($searchText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Text not found: " + $searchText)
    }
    return $rng.Paragraphs(1)
}

$randomText   = "برنامه ای که رنج از اعداد را گرفته و یک عدد تصادفی از بین آنها انتخاب کند."
$shamsiText   = "برنامه که تقویم شمسی را نمایش دهد."
$passwordText = "برنامه ای که یک رمز پویا تولید کند."
$newText      = "برنامه ای که تقویم میلادی را نمایش دهد."

# --- Step 1: swap the text of the "random number" and "Shamsi calendar" bullet paragraphs
# so that, together with a new paragraph inserted before them, the final visual order becomes:
#   [new] calendar(Gregorian), calendar(Shamsi), random-number, dynamic-password, guess-number
$pRandom = Get-ParaByText($randomText)
$pShamsi = Get-ParaByText($shamsiText)

# Edit the later-positioned paragraph first so the earlier paragraph's Start/End stay valid.
if ($pRandom.Range.Start -lt $pShamsi.Range.Start) {
    $shamRange = $d.Range($pShamsi.Range.Start, $pShamsi.Range.End - 1)
    $shamRange.Text = $randomText

    $randRange = $d.Range($pRandom.Range.Start, $pRandom.Range.End - 1)
    $randRange.Text = $shamsiText
} else {
    $randRange = $d.Range($pRandom.Range.Start, $pRandom.Range.End - 1)
    $randRange.Text = $shamsiText

    $shamRange = $d.Range($pShamsi.Range.Start, $pShamsi.Range.End - 1)
    $shamRange.Text = $randomText
}

# --- Step 2: insert a brand new bullet paragraph right before the paragraph that now
# holds the "Shamsi calendar" text (i.e. right before what used to be the "random number" bullet).
$pTarget = Get-ParaByText($shamsiText)
$pTarget.Range.InsertParagraphBefore()

# Re-locate both paragraphs after the insertion shifted indices.
$count = $d.Paragraphs.Count
$newIdx = 0
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq [char]13) {
        $nextPara = $d.Paragraphs($i + 1)
        if ($nextPara.Range.Text.StartsWith($shamsiText)) {
            $newIdx = $i
            break
        }
    }
}
if ($newIdx -eq 0) {
    throw "Could not locate the newly inserted empty paragraph"
}

$newPara = $d.Paragraphs($newIdx)
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newRange.Text = $newText

# Formatting for the new paragraph: spacing 254 auto, green font color 00B050 (keep bold / fonts inherited).
$newPara.LineSpacingRule = 5      # wdLineSpaceMultiple
$newPara.LineSpacing = 12.7       # 254/20 pt
$newRange.Font.Color = 5287936    # RGB(0x00, 0xB0, 0x50) -> 0x00 + 0xB0*256 + 0x50*65536

# --- Step 3: move the "_GoBack" bookmark so it sits right after the newly typed sentence,
# matching where Word would leave it after the last edit was made there.
try {
    $oldBm = $d.Bookmarks("_GoBack")
    $oldBm.Delete()
} catch {
    # bookmark might not exist / already removed - ignore
}
$goBackRange = $d.Range($newRange.End, $newRange.End)
$d.Bookmarks.Add("_GoBack", $goBackRange)

Write-Output "Edit complete"
